# Bitacora de Evaluacion Continua - actualizacion 9 de enero de 2024
# - Captura las notas de "ANDRADE DELGADO BRYANA" en la hoja Concentrando
# - Filtra la hoja Concentrando para mostrar solo esa alumna
# - Actualiza la calificacion de evaluacion continua en la hoja Calificaciones

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Hoja "Concentrando": capturar las puntuaciones de la fila 3
# (ANDRADE DELGADO BRYANA) que antes estaban en 0.
# K3/L3 son formulas y se recalculan solas.
# ---------------------------------------------------------------------
$wsConcentrando = $wb.Worksheets.Item("Concentrando")

$wsConcentrando.Range("F3").Value = 7
$wsConcentrando.Range("G3").Value = 3
$wsConcentrando.Range("H3").Value = 5
$wsConcentrando.Range("I3").Value = 3

# Aplicar AutoFiltro sobre la tabla para dejar visible unicamente a la
# alumna "ANDRADE DELGADO BRYANA" (columna B = campo 2 del rango).
$wsConcentrando.Range("A2:L31").AutoFilter(2, @("ANDRADE DELGADO BRYANA"), 7)

# Ajustar la vista de la hoja: zoom al 120% y celda activa I3.
$wsConcentrando.Activate()
$excel.ActiveWindow.Zoom = 120
$wsConcentrando.Range("I3").Select()

# ---------------------------------------------------------------------
# Hoja "Calificaciones": actualizar la calificacion de evaluacion
# continua (columna C) de la primera alumna. E2 es formula y se
# recalcula sola.
# ---------------------------------------------------------------------
$wsCalificaciones = $wb.Worksheets.Item("Calificaciones")

$wsCalificaciones.Range("C2").Value = 7.83

# Dejar esta hoja como la hoja activa (tal como estaba originalmente) y
# su celda activa en E2, moviendo la vista al tope (B1).
$wsCalificaciones.Activate()
$wsCalificaciones.Range("B1").Select()
$wsCalificaciones.Range("E2").Select()

Write-Output "Bitacora actualizada al 9 de enero de 2024"
